$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_5a_Indikatoren")

# Column D width change (target stored width 39.40234375 chars; engine quantizes
# ColumnWidth writes to 1/7-character steps, so 271/7 = 38.714285714285715 is the
# input that lands closest on the stored value, 39.42857142857143)
$ws.Range("D1").EntireColumn.ColumnWidth = 38.714285714285715

# Row 44
$ws.Range("J44").Value = "Flächendeckender Aufbau bis 2025"
$ws.Range("K44").Value = "Universal Roll-out by 2025"

# Row 46
$ws.Range("D46").Value = "Gini-Koeffizient des Einkommens nach Sozialtransfer"
$ws.Range("E46").Value = "Gini coefficient of income after social transfers"
$ws.Range("J46").Value = "Bis 2030 unterhalb des EU-27-Wertes halten"
$ws.Range("K46").Value = "To be below the EU-27 figure by 2030"

# Row 47
$ws.Range("J47").Value = "Senkung auf unter 30 ha pro Tag bis 2030"

# Row 50
$ws.Range("D50").Value = "Endenergieverbrauch im Güterverkehr"

# Row 51
$ws.Range("D51").Value = "Endenergieverbrauch im Personenverkehr"

# Row 56
$ws.Range("D56").Value = "ba) Rohstoffeinsatz"
$ws.Range("E56").Value = "ba) Use of raw materials"
$ws.Range("L56").Value = "Direkter und indirrekter Rohstoffeinsatz"
$ws.Range("M56").Value = "Direct and indirect use of raw materials"

# Row 57
$ws.Range("D57").Value = "bb) Energieverbrauch"
$ws.Range("E57").Value = "bb) Energy consumption"
$ws.Range("L57").Value = "Direkter und indirrekter Energieverbrauch"
$ws.Range("M57").Value = "Direct and indirect energy consumption"

# Row 58
$ws.Range("D58").Value = "bc) CO2-Emissionen"
$ws.Range("E58").Value = "bc) CO2 emissions"
$ws.Range("L58").Value = "Direkte und indirrekte CO2-Emissionen"
$ws.Range("M58").Value = "Direct and indirect CO2 emissions"

# Row 64
$ws.Range("E64").Value = "aa) Baltic Sea"

# Row 65
$ws.Range("L65").Value = "Stickstoffeinträge in Küsten- und Meeresgewässer der Nordsee"
